# Apply updated cryptocurrency price/volume data per the Nov 16 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.269.60"
$ws.Range("E2").Value = "  -3.10%  "
$ws.Range("D3").Value = "1.983.64"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'244.55"
$ws.Range("E5").Value = "  -3.58%  "
$ws.Range("D6").Value = "'0.629"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D7").Value = "'62.77"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "'56.60"
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("D11").Value = "'0.0801"
$ws.Range("E11").Value = "  +6.45%  "
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").Value = "'0.867"
$ws.Range("E13").Value = "  -4.36%  "
$ws.Range("D14").Value = "'22.51"
$ws.Range("E14").Value = "  +10.71%  "
$ws.Range("D15").Value = "'14.06"
$ws.Range("E15").Value = "  -6.76%  "
$ws.Range("D16").Value = "2.277.43"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("D18").Value = "1.993.66"
$ws.Range("E18").Value = "  -1.83%  "
$ws.Range("D19").Value = "36.084.24"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("D21").Value = "0.0₃0873"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'237.97"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.28"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  -10.16%  "
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'9.78"
$ws.Range("E27").Value = "  +2.80%  "
$ws.Range("D28").Value = "'0.137"
$ws.Range("E28").Value = "  +20.50%  "
$ws.Range("D29").Value = "'159.70"
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("D30").Value = "'19.88"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("E32").Value = "  -4.60%  "
$ws.Range("E33").Value = "  -5.90%  "
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").Value = "'4.39"
$ws.Range("E35").Value = "  -6.04%  "
$ws.Range("D36").Value = "'6.35"
$ws.Range("E36").Value = "  +5.91%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.28"
$ws.Range("E38").Value = "  -7.01%  "
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("E40").Value = "  +15.50%  "
$ws.Range("D41").Value = "'0.0992"
$ws.Range("E41").Value = "  -5.81%  "
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "'2.84"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("E45").Value = "  -3.22%  "
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").Value = "'16.22"
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("D48").Value = "'7.54"
$ws.Range("E48").Value = "  -6.85%  "
$ws.Range("D49").Value = "1.355.71"
$ws.Range("E49").Value = "  -5.09%  "
$ws.Range("D50").Value = "'2.86"
$ws.Range("E50").Value = "  -2.80%  "
$ws.Range("D51").Value = "2.169.51"
$ws.Range("E51").Value = "  -2.15%  "
